$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.03919062573893041
$ws.Range("J2").Value = 0.0391906257389304
$ws.Range("M2").Value = 1.139366
$ws.Range("N2").Value = 3.418098
$ws.Range("O2").Value = 0.1546713947032042
$ws.Range("P2").Value = 0.1546713947032042
$ws.Range("Q2").Value = 0.15428155006
$ws.Range("R2").Value = 1.38853395054
$ws.Range("S2").Value = 0.006061668742331657
$ws.Range("T2").Value = 0.006061668742331656
$ws.Range("I3").Value = 0.03919062573893041
$ws.Range("J3").Value = 0.0391906257389304
$ws.Range("O3").Value = 0.7676983257595695
$ws.Range("P3").Value = 0.7676983257595695
$ws.Range("S3").Value = 0.03008657776524677
$ws.Range("T3").Value = 0.03008657776524676
$ws.Range("I4").Value = 0.03919062573893041
$ws.Range("J4").Value = 0.0391906257389304
$ws.Range("M4").Value = 0.4872916666666667
$ws.Range("N4").Value = 1.461875
$ws.Range("O4").Value = 0.06615089594615092
$ws.Range("P4").Value = 0.06615089594615092
$ws.Range("Q4").Value = 0.06598416458333334
$ws.Range("R4").Value = 0.5938574812499999
$ws.Range("S4").Value = 0.00259249500532053
$ws.Range("T4").Value = 0.002592495005320529
$ws.Range("I5").Value = 0.03919062573893041
$ws.Range("J5").Value = 0.0391906257389304
$ws.Range("M5").Value = 0.08456133333333334
$ws.Range("N5").Value = 0.253684
$ws.Range("O5").Value = 0.0114793835910754
$ws.Range("P5").Value = 0.0114793835910754
$ws.Range("Q5").Value = 0.01145045014666667
$ws.Range("R5").Value = 0.10305405132
$ws.Range("S5").Value = 0.000449884226031455
$ws.Range("T5").Value = 0.000449884226031455
$ws.Range("G6").Value = 3.319753
$ws.Range("H6").Value = 9.959258999999999
$ws.Range("I6").Value = 0.9608093742610696
$ws.Range("J6").Value = 0.9608093742610695
$ws.Range("M6").Value = 1.139366
$ws.Range("N6").Value = 3.418098
$ws.Range("O6").Value = 0.1546713947032042
$ws.Range("P6").Value = 0.1546713947032042
$ws.Range("Q6").Value = 3.782413696598
$ws.Range("R6").Value = 34.041723269382
$ws.Range("S6").Value = 0.1486097259608725
$ws.Range("T6").Value = 0.1486097259608725
$ws.Range("G7").Value = 3.319753
$ws.Range("H7").Value = 9.959258999999999
$ws.Range("I7").Value = 0.9608093742610696
$ws.Range("J7").Value = 0.9608093742610695
$ws.Range("O7").Value = 0.7676983257595695
$ws.Range("P7").Value = 0.7676983257595695
$ws.Range("Q7").Value = 18.77368900552234
$ws.Range("R7").Value = 168.963201049701
$ws.Range("S7").Value = 0.7376117479943228
$ws.Range("T7").Value = 0.7376117479943227
$ws.Range("G8").Value = 3.319753
$ws.Range("H8").Value = 9.959258999999999
$ws.Range("I8").Value = 0.9608093742610696
$ws.Range("J8").Value = 0.9608093742610695
$ws.Range("M8").Value = 0.4872916666666667
$ws.Range("N8").Value = 1.461875
$ws.Range("O8").Value = 0.06615089594615092
$ws.Range("P8").Value = 0.06615089594615092
$ws.Range("Q8").Value = 1.617687972291667
$ws.Range("R8").Value = 14.559191750625
$ws.Range("S8").Value = 0.0635584009408304
$ws.Range("T8").Value = 0.06355840094083039
$ws.Range("G9").Value = 3.319753
$ws.Range("H9").Value = 9.959258999999999
$ws.Range("I9").Value = 0.9608093742610696
$ws.Range("J9").Value = 0.9608093742610695
$ws.Range("M9").Value = 0.08456133333333334
$ws.Range("N9").Value = 0.253684
$ws.Range("O9").Value = 0.0114793835910754
$ws.Range("P9").Value = 0.0114793835910754
$ws.Range("Q9").Value = 0.2807227400173333
$ws.Range("R9").Value = 2.526504660156
$ws.Range("S9").Value = 0.01102949936504394
$ws.Range("T9").Value = 0.01102949936504394
